$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before A - this shifts the existing A:J data to B:K
#    (values, styles and column widths all move with it).
$ws.Columns("A").Insert()

# 2. Fill in the new "_MasterItemID" column.
$ws.Range("A1").Value = "_MasterItemID"
$ws.Range("A2").Value = "D1"
$ws.Range("A3").Value = "M1"
$ws.Range("A4").Value = "D2"
$ws.Range("A5").Value = "D3"

# Header cell needs to be bold like the rest of row 1.
$ws.Range("A1").Font.Bold = $true

# 3. Column A width (closest the host engine can represent to the authored 18.7109375 chars).
$ws.Columns("A:A").ColumnWidth = 17.8333333333333

# 4. Row 4 had two trailing cells (Dim1, Dim2) that got swapped around during the
#    original edit (Dim2 then Dim1) - fix the order post column-insert shift.
$ws.Range("G4").Value = "Dim2"
$ws.Range("H4").Value = "Dim1"

# 5. New row 5 of data ("Super Drill Down" / D3 item).
$ws.Range("B5").Value = "Dimension"
$ws.Range("C5").Value = "Super Drill Down"
$ws.Range("D5").Value = "This drills down through 5 levels"
$ws.Range("F5").Value = "dim;super;tags"
$ws.Range("G5").Value = "Dim1"
$ws.Range("H5").Value = "Dim2"
$ws.Range("I5").Value = "Dim3"
$ws.Range("J5").Value = "AsciiAlpha"
$ws.Range("K5").Value = "AsciiNum"

# 6. Selection follows the same value it used to (moved from old E4 to new G4).
[void]$ws.Range("G4").Select()
